# "Generate Report for Handback"
#
# The handback transform failed for this item, so the localization-status
# report needs to reflect the failure:
#   - the per-language "Status" cells flip from "Ready for handoff" to
#     "Handback transform failed" (Overview!E2/F2, zh-cn!C2, de-de!C2 all
#     share the same text)
#   - the "Error Detail" column on the per-language sheets gets the
#     translationStateItem-not-found message
#   - a couple of columns are widened so the longer text is readable

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$statusText = "Handback transform failed"
$errorText = "The translationStateItem 40c206b14f66a769f1f660de11f1e570d7b3b1a3 is not found."

# Status text, everywhere it shows up.
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$zhcn.Range("C2").Value = $statusText
$dede.Range("C2").Value = $statusText

# New error detail on both per-language sheets.
$zhcn.Range("P2").Value = $errorText
$dede.Range("P2").Value = $errorText

# Widen the columns that now hold longer text.
$overview.Columns.Item(5).ColumnWidth = 23.8
$overview.Columns.Item(6).ColumnWidth = 23.8
$zhcn.Columns.Item(3).ColumnWidth = 23.8
$dede.Columns.Item(3).ColumnWidth = 23.8
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
